$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "ASM_Wt__Q4" coded response wording first (row 25 before the insert)
$ws.Cells.Item(25, 4).Value = "1=Q1`n2=Q2`n3=Q3`n4=Q4"

# Insert a new row at 13 (pushes HT/DM/DysL_/... etc. down by one)
$ws.Rows(13).Insert()

# Fill in the new row 13: MS_5cri (Metabolic syndrome, number of criteria)
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "MS_5cri"
$ws.Cells.Item(13, 3).Value = "Metabolic syndrome, number of criteria"
$ws.Cells.Item(13, 4).Value = "0=0`n1=1`n2=2`n3=3`n4=4`n5=5"
$ws.Cells.Item(13, 5).Value = "Numeric"
$ws.Cells.Item(13, 4).WrapText = $true
$ws.Rows(13).RowHeight = 87

# Renumber the "Number" column for all rows pushed down by the insert
for ($r = 14; $r -le 35; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Update view/selection state to match
$ws.Range("A2:A35").Select()
